$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "October 18th, 2022"
$ws.Range("B8").Value = "URI"
$ws.Range("C8").Value = "UPenn"
$ws.Range("D8").Value = "Barott blue dry shipper"
$ws.Range("E8").Value = "5 Acropora sperm samples in shield"
$ws.Range("F8").Value = "Emma Strand"
$ws.Range("G8").Value = "Ben Glass"
$ws.Range("H8").Value = "NA"

$ws.Range("E8").WrapText = $true
$ws.Range("H8").WrapText = $true
$ws.Range("A8:H8").RowHeight = 17

$ws.Range("H18").Select()
